# Addition of PDMS contaminant series to database
# - Notes sheet: log entry documenting the change
# - Data sheet: 22 new PDMS oligomer contaminant rows (PDMS6..PDMS27)

$wb = $excel.ActiveWorkbook
$dataSheet  = $wb.Worksheets.Item(1)   # LOBSTAHS_rt.windows
$notesSheet = $wb.Worksheets.Item(2)   # Notes

# --- Notes sheet: append change-log row (mirrors the most recent entry's formatting) ---
$notesSheet.Range("A25:C25").Copy($notesSheet.Range("A26:C26"))
$notesSheet.Range("A26").Value = 42779
$notesSheet.Range("B26").Value = "Added PDMS"
$notesSheet.Range("C26").Value = "JEH"

# --- Data sheet: append PDMS6 .. PDMS27 rows with rt window 30 +/- 5 ---
$pdmsNames = @(
    "PDMS6","PDMS7","PDMS8","PDMS9","PDMS10","PDMS11","PDMS12","PDMS13",
    "PDMS14","PDMS15","PDMS16","PDMS17","PDMS18","PDMS19","PDMS20","PDMS21",
    "PDMS22","PDMS23","PDMS24","PDMS25","PDMS26","PDMS27"
)

$row = 73
foreach ($name in $pdmsNames) {
    $dataSheet.Range("A$row").Value = $name
    $dataSheet.Range("A$row").HorizontalAlignment = -4152   # xlRight
    $dataSheet.Range("B$row").Value = 30
    $dataSheet.Range("C$row").Value = 5
    $row = $row + 1
}

# --- Restore a sensible active selection on each sheet (matches end-user state) ---
[void]$dataSheet.Range("F82").Select()
[void]$notesSheet.Range("B27").Select()
[void]$dataSheet.Activate()
